$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new value would otherwise
# be auto-coerced to a number by Excel (these columns are stored as literal text).
$textCells = @('D5','D6','D12','D15','D19','D20','D21','D22','D23','D25','D27','D29','D31','D32','D35','D38','D39','D40','D44','D49','D50')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range('D2').Value = '57.890.92'
$ws.Range('E2').Value = '  +0.96%  '
$ws.Range('D3').Value = '3.104.01'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '526.26'
$ws.Range('E5').Value = '  +1.93%  '
$ws.Range('D6').Value = '141.82'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.106.76'
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '0.384'
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('D13').Value = '3.641.68'
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = '26.29'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('D17').Value = '57.999.31'
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').Value = '3.111.10'
$ws.Range('E18').Value = '  +1.50%  '
$ws.Range('D19').Value = '6.08'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '12.84'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('D21').Value = '8.07'
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('D22').Value = '338.04'
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('D25').Value = '66.26'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').Value = '6.61'
$ws.Range('E29').Value = '  +3.95%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = '7.25'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').Value = '1.86'
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('E33').Value = '  +3.02%  '
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').Value = '154.05'
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  +3.02%  '
$ws.Range('D38').Value = '27.06'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = '1.29'
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('D40').Value = '0.0664'
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = '3.151.77'
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('E42').Value = '  +3.69%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = '36.81'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('E46').Value = '  +6.42%  '
$ws.Range('D47').Value = '2.297.73'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = '0.972'
$ws.Range('E49').Value = '  +4.05%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '20.69'
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('E51').Value = '  +1.99%  '
